$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.467.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.504.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.20%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.502.97"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("E10").Value = "  -7.37%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.962.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.360.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.04%  "
$ws.Range("E16").Value = "  -6.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.508.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.78%  "
$ws.Range("E22").Value = "  -4.20%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.79%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.632.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0903"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "462.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.57%  "
$ws.Range("E34").Value = "  -2.13%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "154.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("E41").Value = "  -3.11%  "
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("E43").Value = "  -6.92%  "
$ws.Range("E44").Value = "  -14.32%  "
$ws.Range("E45").Value = "  -10.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.529"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("E51").Value = "  -2.03%  "
